$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (trial_id) ---

# Rows 52-61 previously held the trial_id as zero-padded text ("04","41",...).
# The random-search results were re-confirmed/re-run, so these become plain numbers
# (same numeric identity, e.g. "04" -> 4) while keeping the rest of the row intact.
$idCol_52_61 = @(4, 41, 11, 35, 16, 3, 42, 5, 23, 21)
for ($i = 0; $i -lt $idCol_52_61.Length; $i++) {
    $ws.Cells.Item(52 + $i, 1).Value = $idCol_52_61[$i]
}

# Rows 62-71 are brand-new rows appended below the previous last row (61).
# Their trial_id is again zero-padded text, so force text via a leading quote
# (the standard Excel way to keep a numeric-looking string, e.g. "09", as text).
$idCol_62_71 = @("09", "13", "37", "14", "36", "20", "18", "30", "03", "27")
for ($i = 0; $i -lt $idCol_62_71.Length; $i++) {
    $ws.Cells.Item(62 + $i, 1).Formula = "'" + $idCol_62_71[$i]
}

# --- Columns B:J (score, l2_lambda, dropout_rate, learning_rate, gru_3, gru_2, gru_1, epocas, numero_anterior) ---
# for rows 52-71.
$rowsBJ = @(
    @(0.03901611516873042, 0.003, 0.1,  0.003, 288, 160, 96,  100, 10),
    @(0.03837998335560163, 0.003, 0.05, 0.001, 288, 128, 64,  100, 10),
    @(0.03795589506626129, 0.002, 0.1,  0.001, 352, 224, 96,  100, 10),
    @(0.03774385154247284, 0.002, 0.05, 0.001, 512, 224, 64,  100, 10),
    @(0.03753180553515752, 0.001, 0.05, 0.003, 448, 160, 64,  100, 10),
    @(0.03710771972934405, 0.002, 0.1,  0.001, 320, 256, 64,  100, 10),
    @(0.03710771848758062, 0.002, 0.05, 0.001, 416, 192, 128, 100, 10),
    @(0.03689567372202873, 0.002, 0.05, 0.001, 256, 160, 96,  100, 10),
    @(0.03647158667445183, 0.003, 0.05, 0.001, 512, 224, 96,  100, 10),
    @(0.03625954315066338, 0.001, 0.05, 0.003, 256, 128, 96,  100, 10),
    @(0.03795589506626129, 0.002, 0.1,  0.001, 416, 192, 128, 100, 9),
    @(0.03774385154247284, 0.002, 0.05, 0.001, 384, 256, 128, 100, 9),
    @(0.03774385154247284, 0.001, 0.1,  0.003, 416, 224, 128, 100, 9),
    @(0.03710771848758062, 0.002, 0.1,  0.003, 480, 256, 96,  100, 9),
    @(0.03710771848758062, 0.003, 0.05, 0.001, 512, 192, 64,  100, 9),
    @(0.03689567496379217, 0.002, 0.05, 0.001, 256, 224, 64,  100, 9),
    @(0.03689567496379217, 0.002, 0.1,  0.001, 480, 256, 128, 100, 9),
    @(0.03668362895647685, 0.001, 0.05, 0.001, 384, 224, 96,  100, 9),
    @(0.03647158667445183, 0.001, 0.05, 0.001, 512, 128, 96,  100, 9),
    @(0.03583545486132304, 0.001, 0.1,  0.001, 384, 192, 128, 100, 9)
)

$data = New-Object 'object[,]' $rowsBJ.Length, 9
for ($i = 0; $i -lt $rowsBJ.Length; $i++) {
    $row = $rowsBJ[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $data[$i, $j] = $row[$j]
    }
}
$ws.Range("B52:J71").Value = $data
